# Generate Report for Handoff
# b.md is ready for handoff: update Status on Overview/zh-cn/de-de sheets,
# and record the new Latest Handoff File / Datetime for zh-cn + de-de,
# including updating the displayed hyperlink text for the new handoff file.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: Status column for b.md (row 3) ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B3").Value = "Ready for handoff"
$ovw.Range("C3").Value = "Ready for handoff"

# ---- zh-cn sheet: Status + Latest Handoff File/Datetime for b.md (row 3) ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("D3").Value = "2016-01-27 07:50:23"

# Rebuild the sheet's hyperlinks in-place so the only visible change is the
# updated display text for C3 (the rest keep their original target + text).
$zhLinks = @(
    @{ Ref = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/3c676b4e880d569ce4903684c185d7277b069def/e2e/a.md"; Display = "a.md" },
    @{ Ref = "C2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09126ad7209111689aad3317371195111403b551/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "E2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/55293bf0fad742a0e147395bb3afabf64b63a421/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/65ea3a31312b3dee3ca49377ef3f4f1b5c11c6a2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/3c676b4e880d569ce4903684c185d7277b069def/e2e/b.md"; Display = "b.md" },
    @{ Ref = "C3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09126ad7209111689aad3317371195111403b551/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" },
    @{ Ref = "E3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/55293bf0fad742a0e147395bb3afabf64b63a421/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/65ea3a31312b3dee3ca49377ef3f4f1b5c11c6a2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "A4"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/3c676b4e880d569ce4903684c185d7277b069def/.localization-config"; Display = ".localization-config" }
)
$zh.Hyperlinks.Delete()
foreach ($link in $zhLinks) {
    $zh.Hyperlinks.Add($zh.Range($link.Ref), $link.Address, "", "", $link.Display)
}

# ---- de-de sheet: Status + Latest Handoff File/Datetime for b.md (row 3) ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("D3").Value = "2016-01-27 07:50:34"

$deLinks = @(
    @{ Ref = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/3c676b4e880d569ce4903684c185d7277b069def/e2e/a.md"; Display = "a.md" },
    @{ Ref = "C2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b823569a9fb1234b28604621a0ff05840a23c0db/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "E2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/dc5e1bf716c34ed6364344cc407b41e653da976a/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3f7afb6143d0e80067ac570e00090c0d8b54316c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/tianzh/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/3c676b4e880d569ce4903684c185d7277b069def/e2e/b.md"; Display = "b.md" },
    @{ Ref = "C3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b823569a9fb1234b28604621a0ff05840a23c0db/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" },
    @{ Ref = "E3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/dc5e1bf716c34ed6364344cc407b41e653da976a/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3f7afb6143d0e80067ac570e00090c0d8b54316c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/tianzh/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "A4"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/3c676b4e880d569ce4903684c185d7277b069def/.localization-config"; Display = ".localization-config" }
)
$de.Hyperlinks.Delete()
foreach ($link in $deLinks) {
    $de.Hyperlinks.Add($de.Range($link.Ref), $link.Address, "", "", $link.Display)
}
